# Workbook already open; grab it and the sheet that holds the tracker data
# (the "Sheet1" tab -- it is already the active sheet/tab in this workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the next five challenges as passed by replacing the "Not Passed..."
# text in D81:D85 with the corresponding "Passed..." text. Setting .Value to
# a brand-new string causes Excel to append new shared-string entries.
$ws.Range("D81").Value = "PassedCounting Cards"
$ws.Range("D82").Value = "PassedBuild JavaScript Objects"
$ws.Range("D83").Value = "PassedAccessing Object Properties with Dot Notation"
$ws.Range("D84").Value = "PassedAccessing Object Properties with Bracket Notation"
$ws.Range("D85").Value = "PassedAccessing Object Properties with Variables"

# Move the view/selection to the rows that were just edited, matching the
# author's on-screen state when the change was saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 78
$win.ScrollColumn = 1
$ws.Range("D81:D85").Select() | Out-Null
